$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New protein identifiers / gene names to append below existing list
$ws.Range("A4").Value = "SPCS_HUMAN"
$ws.Range("A5").Value = "SEPP1_HUMAN"
$ws.Range("A6").Value = "ISCU_HUMAN"
$ws.Range("A7").Value = "Q92911"
$ws.Range("A8").Value = "P52789"
$ws.Range("A9").Value = "Q9UPP1"
$ws.Range("A10").Value = "O43772"
$ws.Range("A11").Value = "P21796"

# Emphasize A6 (ISCU_HUMAN) with bold, larger font
$ws.Range("A6").Font.Bold = $true
$ws.Range("A6").Font.Size = 12
$ws.Rows(6).RowHeight = 15.75

# Update current selection to mirror the authored state
$ws.Range("D10").Select()
